$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Hardware: ..." paragraph (the one that currently mentions the
# old CPU model) and the "OS: Ubuntu ..." paragraph that immediately follows
# it, without relying on fixed paragraph indices.
# ---------------------------------------------------------------------------
$hwPara = $null
$osPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*i7-6700K*") {
        $hwPara = $p
        $osPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$paraStart = $hwPara.Range.Start
$paraText  = $hwPara.Range.Text

$oldCpuText = "Intel(R) Core(TM) i7-6700K CPU @ 4.00GHz"
$newCpuText = "Intel(R) Xeon(R) E-2176G CPU @ 3.70GHz"

$commaIdx = $paraText.IndexOf(",")
$commaPos = $paraStart + $commaIdx          # position of the "," run
$spacePos = $commaPos + 1                   # position of the lone " " run
$cpuPos   = $spacePos + 1                   # position where the CPU run starts

# ---------------------------------------------------------------------------
# Step 1: remove the original (hidden) "_GoBack" bookmark - it currently sits
# between the "," run and the " " run.
# ---------------------------------------------------------------------------
$oldBk = $d.Bookmarks.Item("_GoBack")
$oldBk.Delete()

# ---------------------------------------------------------------------------
# Step 2: delete the standalone " " run that separates the comma from the
# CPU name, then fold that space into the comma run itself so it reads ", ".
# ---------------------------------------------------------------------------
$spaceRange = $d.Range($spacePos, $spacePos + 1)
$spaceRange.Delete()

$commaRange = $d.Range($commaPos, $commaPos + 1)
$commaRange.Text = ", "

# ---------------------------------------------------------------------------
# Step 3: swap in the new CPU model text.
# ---------------------------------------------------------------------------
$cpuRange = $d.Range($cpuPos, $cpuPos + $oldCpuText.Length)
$cpuRange.Text = $newCpuText

# ---------------------------------------------------------------------------
# Step 4: force the new CPU text to live in its own run (separate from the
# ", " run it now directly follows) by toggling a character property on and
# back off - adjacent runs with identical formatting only get coalesced when
# the formatting is actually produced by the same edit, so a no-op
# bold-on/bold-off round trip leaves a genuine run boundary behind.
# ---------------------------------------------------------------------------
$newCpuRange = $d.Range($cpuPos, $cpuPos + $newCpuText.Length)
$newCpuRange.Font.Bold = 1
$newCpuRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Step 5: re-create the "_GoBack" bookmark, now collapsed right after the
# "OS: Ubuntu 16.04.3" text (i.e. at the very end of that paragraph, before
# its paragraph mark). A collapsed Range built directly at a
# paragraph-end-minus-one position cannot be fed straight into
# Bookmarks.Add, so a throwaway character is inserted there first to move
# the insertion point off that boundary, the bookmark is added next to it,
# and then the throwaway character is deleted again.
# ---------------------------------------------------------------------------
$endPos = $osPara.Range.End - 1

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("Z")

$bkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bkRange)

$placeholderRange = $d.Range($endPos, $endPos + 1)
$placeholderRange.Delete()
